# TC13_Trials_Filter_Diagnosis-GallBladder: add the Cypher "Gall bladder"
# query to the startup sheet (A2), matching sharedStrings + row layout,
# and leave the sheet's selection on the B2:B5 block.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Add the Cypher query text to A2 (wraps; row height grows to match the authored value)
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Gall bladder carcinoma (adenocarcinoma)''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

$ws.Rows.Item(2).RowHeight = 101.5

# Update the sheet's selection/view state
$ws.Range("B2:B5").Select()
